$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Förändrad" (Changed) column C holds a date serial that was bumped
# by one day (2023-09-09 -> 2023-09-10) for every data row (rows 2-484).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 484 }

$range = $ws.Range("C2:C$lastRow")
$range.Value = 45179
